$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows 5 and 6 (Mice, Rats), shifting remaining rows up
$ws.Range("A5:A6").EntireRow.Delete()

# Update selection to match the resulting state
$ws.Range("A5:XFD6").Select()
